$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 25: merge "Invoicing" + "Refund Services" info into a single "Invoicing & Refund Services" row
$ws.Range("A25").Value = "Invoicing & Refund Services "
$ws.Range("B25").Value = "Invoicing / Sales Posting / Credit Note / Debit Note, Payment Ledger, Refund Ledger"

# Row 26: repurpose the now-freed row for a new "Material Management" entry
$ws.Range("A26").Value = "Material Management "
$ws.Range("B26").Value = "Inventory Management and Ledger "
$ws.Range("C26").Value = "SAP IM "

# Update view: scroll so row 17 is the top row, and select B30
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("B30").Select()
